$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 64
$ws.Range("H64").Value = 4331
$ws.Range("J64").Value = 4746.5
$ws.Range("L64").Value = 4746.5
$ws.Range("N64").Value = -5242.5
# row 67
$ws.Range("H67").Value = 4331
$ws.Range("J67").Value = 4746.5
$ws.Range("L67").Value = 4746.5
$ws.Range("N67").Value = -6462.5
# row 101
$ws.Range("H101").Value = 14287151
$ws.Range("J101").Value = 1883.4
$ws.Range("L101").Value = 5650.200000000001
$ws.Range("N101").Value = -8894.200000000001
# row 112
$ws.Range("H112").Value = 2507.25
$ws.Range("I112").Value = 1014
$ws.Range("J112").Value = 3005
$ws.Range("K112").Value = 3042
$ws.Range("L112").Value = 9015
$ws.Range("N112").Value = -11231
$ws.Range("M112").Value = -1934
# row 135
$ws.Range("H135").Value = 1221.1818
$ws.Range("I135").Value = 1254.5
$ws.Range("K135").Value = 11290.5
$ws.Range("M135").Value = -8755.5
# row 137
$ws.Range("H137").Value = 1899.5
$ws.Range("I137").Value = 1666.5
$ws.Range("J137").Value = 2598.5
$ws.Range("K137").Value = 4999.5
$ws.Range("L137").Value = 7795.5
$ws.Range("M137").Value = -2449.5
$ws.Range("N137").Value = -12895.5
# row 138
$ws.Range("H138").Value = 3500.6538
$ws.Range("J138").Value = 3417.3906
$ws.Range("L138").Value = 10252.1718
$ws.Range("N138").Value = -20532.1718

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 2
$ws.Range("H2").Value = 1753.5264
$ws.Range("I2").Value = 811.36365
$ws.Range("K2").Value = 811.36365
$ws.Range("M2").Value = -698.36365
# row 32
$ws.Range("H32").Value = 5005.1514
$ws.Range("I32").Value = 5090
$ws.Range("J32").Value = 3690
$ws.Range("K32").Value = 5090
$ws.Range("L32").Value = 3690
$ws.Range("M32").Value = -4803
$ws.Range("N32").Value = -4264
# row 61
$ws.Range("H61").Value = 5588.6665
$ws.Range("I61").Value = 1973.3334
$ws.Range("K61").Value = 1973.3334
$ws.Range("M61").Value = -1761.3334
# row 116
$ws.Range("H116").Value = 1753.5264
$ws.Range("I116").Value = 811.36365
$ws.Range("K116").Value = 811.36365
$ws.Range("M116").Value = 1482.63635
# row 132
$ws.Range("H132").Value = 2067.8125
$ws.Range("I132").Value = 1871.4
$ws.Range("J132").Value = 5014
$ws.Range("K132").Value = 5614.200000000001
$ws.Range("L132").Value = 15042
$ws.Range("M132").Value = -3084.200000000001
$ws.Range("N132").Value = -20102
# row 136
$ws.Range("H136").Value = 5588.6665
$ws.Range("I136").Value = 1973.3334
$ws.Range("K136").Value = 5920.0002
$ws.Range("M136").Value = -3370.0002

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 3
$ws.Range("H3").Value = 1753.5264
$ws.Range("I3").Value = 811.36365
$ws.Range("K3").Value = 811.36365
$ws.Range("M3").Value = -697.36365
# row 36
$ws.Range("H36").Value = 1170
$ws.Range("I36").Value = 804
$ws.Range("J36").Value = 3000
$ws.Range("K36").Value = 804
$ws.Range("L36").Value = 3000
$ws.Range("M36").Value = -270
$ws.Range("N36").Value = -4068
# row 37
$ws.Range("H37").Value = 5000
$ws.Range("J37").Value = 5000
$ws.Range("L37").Value = 5000
$ws.Range("N37").Value = -5274
# row 134
$ws.Range("H134").Value = 647.5
$ws.Range("I134").Value = 647.5
$ws.Range("K134").Value = 1942.5
$ws.Range("M134").Value = 592.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 58
$ws.Range("H58").Value = 3384.3333
$ws.Range("I58").Value = 3461.2
$ws.Range("K58").Value = 3461.2
$ws.Range("M58").Value = -3258.2
# row 99
$ws.Range("H99").Value = 11387.538
$ws.Range("J99").Value = 5991.3335
$ws.Range("L99").Value = 5991.3335
$ws.Range("N99").Value = -8987.333500000001
# row 107
$ws.Range("H107").Value = 621.0833
$ws.Range("I107").Value = 552.55554
$ws.Range("K107").Value = 552.55554
$ws.Range("M107").Value = 1367.44446
# row 126
$ws.Range("H126").Value = 11387.538
$ws.Range("J126").Value = 5991.3335
$ws.Range("L126").Value = 17974.0005
$ws.Range("N126").Value = -22914.0005
# row 136
$ws.Range("H136").Value = 3384.3333
$ws.Range("I136").Value = 3461.2
$ws.Range("K136").Value = 10383.6
$ws.Range("M136").Value = -7833.599999999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 4
$ws.Range("H4").Value = 27500300
$ws.Range("I4").Value = 27500300
$ws.Range("K4").Value = 82500900
$ws.Range("M4").Value = -82500788
# row 121
$ws.Range("H121").Value = 10008.368
$ws.Range("J121").Value = 4515
$ws.Range("L121").Value = 13545
$ws.Range("N121").Value = -16165
# row 122
$ws.Range("H122").Value = 506.7143
$ws.Range("J122").Value = 503.2
$ws.Range("L122").Value = 4528.8
$ws.Range("N122").Value = -9428.799999999999
# row 127
$ws.Range("H127").Value = 1992.6
$ws.Range("J127").Value = 1992.6
$ws.Range("L127").Value = 5977.799999999999
$ws.Range("N127").Value = -15897.8
# row 131
$ws.Range("H131").Value = 2879.6667
$ws.Range("J131").Value = 3100
$ws.Range("L131").Value = 9300
$ws.Range("N131").Value = -19380
# row 134
$ws.Range("H134").Value = 18389.715
$ws.Range("J134").Value = 18389.715
$ws.Range("L134").Value = 55169.145
$ws.Range("N134").Value = -65309.145

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 12
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
# row 62
$ws.Range("H62").Value = 90077
$ws.Range("I62").Value = 90077
$ws.Range("K62").Value = 90077
$ws.Range("M62").Value = -89391
# row 65
$ws.Range("H65").Value = 90077
$ws.Range("I65").Value = 90077
$ws.Range("K65").Value = 270231
$ws.Range("M65").Value = -266799
# row 70
$ws.Range("H70").Value = 4494.2
$ws.Range("I70").Value = 4535
$ws.Range("J70").Value = 4399
$ws.Range("K70").Value = 4535
$ws.Range("L70").Value = 4399
$ws.Range("M70").Value = -4265
$ws.Range("N70").Value = -4939
# row 73
$ws.Range("H73").Value = 4494.2
$ws.Range("I73").Value = 4535
$ws.Range("J73").Value = 4399
$ws.Range("K73").Value = 4535
$ws.Range("L73").Value = 4399
$ws.Range("M73").Value = -3599
$ws.Range("N73").Value = -6271
# row 97
$ws.Range("H97").Value = 711
$ws.Range("I97").Value = 605.4
$ws.Range("K97").Value = 605.4
$ws.Range("M97").Value = -109.4
# row 132
$ws.Range("H132").Value = 1962
$ws.Range("I132").Value = 1962
$ws.Range("K132").Value = 5886
$ws.Range("M132").Value = -3356
# row 136
$ws.Range("H136").Value = 20994.166
$ws.Range("J136").Value = 20994.166
$ws.Range("L136").Value = 62982.49800000001
$ws.Range("N136").Value = -68082.49800000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 136
$ws.Range("H136").Value = 2962
$ws.Range("I136").Value = 1949.6666
$ws.Range("K136").Value = 5848.9998
$ws.Range("M136").Value = -3298.9998

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 54
$ws.Range("H54").Value = 23657.715
$ws.Range("I54").Value = 14121
$ws.Range("J54").Value = 47499.5
$ws.Range("K54").Value = 14121
$ws.Range("L54").Value = 47499.5
$ws.Range("M54").Value = -13601
$ws.Range("N54").Value = -48539.5
# row 103
$ws.Range("H103").Value = 35650.5
$ws.Range("J103").Value = 35650.5
$ws.Range("L103").Value = 35650.5
$ws.Range("N103").Value = -37994.5
